$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells, matching the style of the existing header row
$ws.Range("E1").Value = "Execution Time (ms)"
$ws.Range("F1").Value = "Memory Usage (B)"
$ws.Range("A1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)

# Execution time (ms) values for each model
$executionTimes = @(3.807600005529821, 3.063299984205514, 15.64580001286231, 4.860599990934134, 1.350300008198246)

# Memory usage (B) values for each model
$memoryUsages = @(0, 0, 0, 0, 0)

for ($i = 0; $i -lt 5; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $executionTimes[$i]
    $ws.Cells.Item($row, 6).Value = $memoryUsages[$i]
}
